$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 9.1131411301072
$ws.Range("E2").Value = 7.264383834843784
$ws.Range("F2").Value = 9.904341455780914
$ws.Range("G2").Value = 7.792725627883438
